# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with pristine default style (no explicit style index) used to restore style
$refStyle = $ws.Range("D18").Style

# --- Column D updates whose new value LOOKS like a number (e.g. '210.68') ---
# Force cell format to Text before assigning so Excel keeps the literal string
# instead of converting it to a floating point number, then restore the original
# (default) cell style so formatting stays identical to before the edit.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.68"
$ws.Range("D5").Style = $refStyle
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5200"
$ws.Range("D6").Style = $refStyle
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2624"
$ws.Range("D8").Style = $refStyle
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07553"
$ws.Range("D11").Style = $refStyle
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.440"
$ws.Range("D13").Style = $refStyle
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5493"
$ws.Range("D14").Style = $refStyle
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008016"
$ws.Range("D15").Style = $refStyle
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.44"
$ws.Range("D16").Style = $refStyle
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.760"
$ws.Range("D19").Style = $refStyle
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.11"
$ws.Range("D20").Style = $refStyle
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.210"
$ws.Range("D22").Style = $refStyle
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "150.07"
$ws.Range("D24").Style = $refStyle
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1240"
$ws.Range("D25").Style = $refStyle
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.500"
$ws.Range("D26").Style = $refStyle
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06312"
$ws.Range("D28").Style = $refStyle
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.351"
$ws.Range("D29").Style = $refStyle
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.524"
$ws.Range("D31").Style = $refStyle
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6063"
$ws.Range("D35").Style = $refStyle
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.759"
$ws.Range("D37").Style = $refStyle
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.120"
$ws.Range("D39").Style = $refStyle
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8636"
$ws.Range("D41").Style = $refStyle
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("D42").Style = $refStyle
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.50"
$ws.Range("D43").Style = $refStyle
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("D45").Style = $refStyle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.56"
$ws.Range("D46").Style = $refStyle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9990"
$ws.Range("D47").Style = $refStyle
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.049"
$ws.Range("D48").Style = $refStyle
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4242"
$ws.Range("D50").Style = $refStyle
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.930"
$ws.Range("D51").Style = $refStyle

# --- Column D updates whose new value is NOT a valid number (e.g. '26.163.31') ---
# These remain text automatically, no special handling required.
$ws.Range("D2").Value = "26.163.31"
$ws.Range("D3").Value = "1.669.76"
$ws.Range("D12").Value = "1.674.01"
$ws.Range("D17").Value = "26.162.46"
$ws.Range("D38").Value = "1.112.92"
$ws.Range("D44").Value = "1.823.42"

# --- Column E (percentage change) updates ---
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("E5").Value = "  -2.60%  "
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("E14").Value = "  -4.32%  "
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("E21").Value = "  -3.86%  "
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("E32").Value = "  -4.30%  "
$ws.Range("E33").Value = "  -2.12%  "
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  +1.45%  "
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  -2.37%  "
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("E51").Value = "  -1.14%  "
